# Auto-generated: apply meteocat daily summary refresh
# (timestamps, percentages, pressures, temperatures, radiation, wind gust)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-24 21:48:18"
$ws.Range("O2").Value = "5.7 °C"
$ws.Range("E3").Value = "2026-02-24 21:48:21"
$ws.Range("E4").Value = "2026-02-24 21:48:23"
$ws.Range("O4").Value = "12.7 °C"
$ws.Range("E5").Value = "2026-02-24 21:48:26"
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = "29%"
$ws.Range("K5").Value = "15.6 MJ/m2"
$ws.Range("E6").Value = "2026-02-24 21:48:28"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "70%"
$ws.Range("J6").Value = "1019.5 hPa"
$ws.Range("O6").Value = "14.0 °C"
$ws.Range("E7").Value = "2026-02-24 21:48:30"
$ws.Range("J7").Value = "1020.1 hPa"
$ws.Range("E8").Value = "2026-02-24 21:48:32"
$ws.Range("J8").Value = "1019.5 hPa"
$ws.Range("E9").Value = "2026-02-24 21:48:35"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "82%"
$ws.Range("E10").Value = "2026-02-24 21:48:37"
$ws.Range("O10").Value = "11.1 °C"
$ws.Range("E11").Value = "2026-02-24 21:48:40"
$ws.Range("O11").Value = "9.1 °C"
$ws.Range("E12").Value = "2026-02-24 21:48:42"
$ws.Range("O12").Value = "10.5 °C"
$ws.Range("E13").Value = "2026-02-24 21:48:44"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "63%"
$ws.Range("O13").Value = "6.7 °C"
$ws.Range("E14").Value = "2026-02-24 21:48:47"
$ws.Range("N14").Value = "6.0 °C 21:29 TU"
$ws.Range("O14").Value = "11.4 °C"
$ws.Range("E15").Value = "2026-02-24 21:48:49"
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "77%"
$ws.Range("O15").Value = "11.8 °C"
$ws.Range("E16").Value = "2026-02-24 21:48:51"
$ws.Range("E17").Value = "2026-02-24 21:48:54"
$ws.Range("E18").Value = "2026-02-24 21:48:56"
$ws.Range("J18").Value = "1020.0 hPa"
$ws.Range("O18").Value = "11.2 °C"
$ws.Range("E19").Value = "2026-02-24 21:48:58"
$ws.Range("E20").Value = "2026-02-24 21:49:00"
$ws.Range("E21").Value = "2026-02-24 21:49:03"
$ws.Range("J21").Value = "1021.8 hPa"
$ws.Range("E22").Value = "2026-02-24 21:49:05"
$ws.Range("H22").NumberFormat = "@"
$ws.Range("H22").Value = "26%"
$ws.Range("E23").Value = "2026-02-24 21:49:07"
$ws.Range("O23").Value = "4.7 °C"
$ws.Range("E24").Value = "2026-02-24 21:49:10"
$ws.Range("L24").Value = "13.3 km/h - 80º 21:09 TU"
$ws.Range("E25").Value = "2026-02-24 21:49:12"
$ws.Range("O25").Value = "6.6 °C"
$ws.Range("E26").Value = "2026-02-24 21:49:15"
$ws.Range("E27").Value = "2026-02-24 21:49:17"
$ws.Range("H27").NumberFormat = "@"
$ws.Range("H27").Value = "33%"
$ws.Range("E28").Value = "2026-02-24 21:49:19"
$ws.Range("J28").Value = "1020.0 hPa"
$ws.Range("O28").Value = "11.5 °C"
$ws.Range("E29").Value = "2026-02-24 21:49:22"
$ws.Range("O29").Value = "10.1 °C"
$ws.Range("E30").Value = "2026-02-24 21:49:24"
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H30").Value = "76%"
$ws.Range("O30").Value = "13.0 °C"
$ws.Range("E31").Value = "2026-02-24 21:49:26"
$ws.Range("J31").Value = "1019.0 hPa"
$ws.Range("E32").Value = "2026-02-24 21:49:29"
$ws.Range("H32").NumberFormat = "@"
$ws.Range("H32").Value = "69%"
$ws.Range("O32").Value = "7.1 °C"
$ws.Range("E33").Value = "2026-02-24 21:49:31"
$ws.Range("E34").Value = "2026-02-24 21:49:33"
$ws.Range("H34").NumberFormat = "@"
$ws.Range("H34").Value = "49%"
$ws.Range("E35").Value = "2026-02-24 21:49:36"
$ws.Range("O35").Value = "13.0 °C"
$ws.Range("E36").Value = "2026-02-24 21:49:38"
$ws.Range("O36").Value = "12.9 °C"
$ws.Range("E37").Value = "2026-02-24 21:49:41"
$ws.Range("H37").NumberFormat = "@"
$ws.Range("H37").Value = "72%"
$ws.Range("O37").Value = "8.6 °C"
$ws.Range("E38").Value = "2026-02-24 21:49:43"
$ws.Range("H38").NumberFormat = "@"
$ws.Range("H38").Value = "73%"
$ws.Range("O38").Value = "11.9 °C"
$ws.Range("E39").Value = "2026-02-24 21:49:46"
$ws.Range("E40").Value = "2026-02-24 21:49:48"
$ws.Range("O40").Value = "8.5 °C"
$ws.Range("E41").Value = "2026-02-24 21:49:50"
$ws.Range("H41").NumberFormat = "@"
$ws.Range("H41").Value = "80%"
$ws.Range("O41").Value = "10.8 °C"
$ws.Range("E42").Value = "2026-02-24 21:49:53"
$ws.Range("O42").Value = "11.2 °C"
$ws.Range("E43").Value = "2026-02-24 21:49:55"
$ws.Range("O43").Value = "10.7 °C"
$ws.Range("E44").Value = "2026-02-24 21:49:57"
$ws.Range("E45").Value = "2026-02-24 21:50:00"
$ws.Range("O45").Value = "10.0 °C"
$ws.Range("E46").Value = "2026-02-24 21:50:02"
